$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 169.86667
$ws.Range("J9").Value = 217.5
$ws.Range("L9").Value = 217.5
$ws.Range("N9").Value = -555.5

$ws.Range("H17").Value = 5278.0454
$ws.Range("J17").Value = 5501.4443
$ws.Range("L17").Value = 16504.3329
$ws.Range("N17").Value = -16840.3329

$ws.Range("H53").Value = 4757.8335
$ws.Range("J53").Value = 9305.916999999999
$ws.Range("L53").Value = 9305.916999999999
$ws.Range("N53").Value = -10579.917

$ws.Range("H86").Value = 2381.2856
$ws.Range("I86").Value = 2266.111
$ws.Range("K86").Value = 2266.111
$ws.Range("M86").Value = -1143.111

$ws.Range("H89").Value = 2381.2856
$ws.Range("I89").Value = 2266.111
$ws.Range("K89").Value = 11330.555
$ws.Range("M89").Value = -5714.555

$ws.Range("H98").Value = 2859.3333
$ws.Range("I98").Value = 2783.5
$ws.Range("J98").Value = 3011
$ws.Range("K98").Value = 2783.5
$ws.Range("L98").Value = 3011
$ws.Range("M98").Value = -1285.5
$ws.Range("N98").Value = -6007

$ws.Range("H122").Value = 2859.3333
$ws.Range("I122").Value = 2783.5
$ws.Range("J122").Value = 3011
$ws.Range("K122").Value = 8350.5
$ws.Range("L122").Value = 9033
$ws.Range("M122").Value = -5900.5
$ws.Range("N122").Value = -13933

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2925633
$ws.Range("I2").Value = 5556883
$ws.Range("K2").Value = 5556883
$ws.Range("M2").Value = -5556770

$ws.Range("H61").Value = 4202.75
$ws.Range("J61").Value = 4361.875
$ws.Range("L61").Value = 4361.875
$ws.Range("N61").Value = -4785.875

$ws.Range("H74").Value = 53037.625
$ws.Range("I74").Value = 7273.1875
$ws.Range("K74").Value = 7273.1875
$ws.Range("M74").Value = -6399.1875

$ws.Range("H77").Value = 53037.625
$ws.Range("I77").Value = 7273.1875
$ws.Range("K77").Value = 36365.9375
$ws.Range("M77").Value = -31997.9375

$ws.Range("H116").Value = 2925633
$ws.Range("I116").Value = 5556883
$ws.Range("K116").Value = 5556883
$ws.Range("M116").Value = -5554589

$ws.Range("H132").Value = 2915.8647
$ws.Range("I132").Value = 2280.2
$ws.Range("K132").Value = 6840.599999999999
$ws.Range("M132").Value = -4310.599999999999

$ws.Range("H136").Value = 4202.75
$ws.Range("J136").Value = 4361.875
$ws.Range("L136").Value = 13085.625
$ws.Range("N136").Value = -18185.625

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2925633
$ws.Range("I3").Value = 5556883
$ws.Range("K3").Value = 5556883
$ws.Range("M3").Value = -5556769

$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = ""
$ws.Range("N76").Value = ""

$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = ""
$ws.Range("N79").Value = ""

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H37").Value = 30000
$ws.Range("J37").Value = 30000
$ws.Range("L37").Value = 30000
$ws.Range("N37").Value = -30214

$ws.Range("H44").Value = 60000
$ws.Range("J44").Value = 30000
$ws.Range("L44").Value = 30000
$ws.Range("N44").Value = -30884

$ws.Range("H45").Value = 14499.5
$ws.Range("J45").Value = 15000
$ws.Range("L45").Value = 15000
$ws.Range("N45").Value = -16186

$ws.Range("H51").Value = 44000
$ws.Range("J51").Value = 44000
$ws.Range("L51").Value = 44000
$ws.Range("N51").Value = -45472

$ws.Range("H61").Value = 44000
$ws.Range("J61").Value = 44000
$ws.Range("L61").Value = 44000
$ws.Range("N61").Value = -44696

$ws.Range("H62").Value = 3000
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 3000
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = ""
$ws.Range("M62").Value = ""
$ws.Range("N62").Value = -4248

$ws.Range("H65").Value = 3000
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 3000
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 15000
$ws.Range("M65").Value = ""
$ws.Range("N65").Value = -21240

$ws.Range("H86").Value = 10816.214
$ws.Range("I86").Value = 9838.25
$ws.Range("J86").Value = 12120.167
$ws.Range("K86").Value = 9838.25
$ws.Range("L86").Value = 12120.167
$ws.Range("M86").Value = -8715.25
$ws.Range("N86").Value = -14366.167

$ws.Range("H89").Value = 10816.214
$ws.Range("I89").Value = 9838.25
$ws.Range("J89").Value = 12120.167
$ws.Range("K89").Value = 49191.25
$ws.Range("L89").Value = 60600.835
$ws.Range("M89").Value = -43575.25
$ws.Range("N89").Value = -71832.83499999999

$ws.Range("H92").Value = 33331
$ws.Range("J92").Value = 33331
$ws.Range("L92").Value = 33331
$ws.Range("N92").Value = -38323

$ws.Range("H97").Value = 21500
$ws.Range("J97").Value = 21500
$ws.Range("L97").Value = 21500
$ws.Range("N97").Value = -23482

$ws.Range("H99").Value = 4337.636
$ws.Range("I99").Value = 3969.8333
$ws.Range("K99").Value = 3969.8333
$ws.Range("M99").Value = -2471.8333

$ws.Range("H102").Value = 44667.75
$ws.Range("J102").Value = 44667.75
$ws.Range("L102").Value = 44667.75
$ws.Range("N102").Value = -49535.75

$ws.Range("H109").Value = 32998
$ws.Range("J109").Value = 32998
$ws.Range("L109").Value = 32998
$ws.Range("N109").Value = -35078

$ws.Range("H126").Value = 4337.636
$ws.Range("I126").Value = 3969.8333
$ws.Range("K126").Value = 11909.4999
$ws.Range("M126").Value = -9439.499899999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 4134411.8
$ws.Range("I4").Value = 6145503
$ws.Range("J4").Value = 112229.75
$ws.Range("K4").Value = 18436509
$ws.Range("L4").Value = 336689.25
$ws.Range("M4").Value = -18436397
$ws.Range("N4").Value = -336913.25

$ws.Range("H60").Value = 1175
$ws.Range("I60").Value = 1175
$ws.Range("K60").Value = 3525
$ws.Range("M60").Value = -3274

$ws.Range("H122").Value = 991.7778
$ws.Range("I122").Value = 865.875
$ws.Range("J122").Value = 1999
$ws.Range("K122").Value = 7792.875
$ws.Range("L122").Value = 17991
$ws.Range("M122").Value = -5342.875
$ws.Range("N122").Value = -22891

$ws.Range("H129").Value = 850.7
$ws.Range("I129").Value = 675.8570999999999
$ws.Range("K129").Value = 2027.5713
$ws.Range("M129").Value = 2972.4287

$ws.Range("H132").Value = 2083.75
$ws.Range("I132").Value = 1290.4546
$ws.Range("K132").Value = 11614.0914
$ws.Range("M132").Value = -9084.091400000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H29").Value = 2250
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").Value = ""

$ws.Range("H107").Value = 1859.125
$ws.Range("I107").Value = 1880.6
$ws.Range("K107").Value = 1880.6
$ws.Range("M107").Value = 39.40000000000009

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3566.95
$ws.Range("I7").Value = 1419.2
$ws.Range("K7").Value = 1419.2
$ws.Range("M7").Value = -1307.2

$ws.Range("H22").Value = 112860.25
$ws.Range("J22").Value = 2165.6667
$ws.Range("L22").Value = 2165.6667
$ws.Range("N22").Value = -2755.6667

$ws.Range("H27").Value = 112860.25
$ws.Range("J27").Value = 2165.6667
$ws.Range("L27").Value = 2165.6667
$ws.Range("N27").Value = -2379.6667

$ws.Range("H40").Value = 5380.5557
$ws.Range("J40").Value = 6659.3335
$ws.Range("L40").Value = 6659.3335
$ws.Range("N40").Value = -6931.3335

$ws.Range("H122").Value = 5457.136
$ws.Range("J122").Value = 8683.286
$ws.Range("L122").Value = 26049.858
$ws.Range("N122").Value = -30949.858

$ws.Range("H126").Value = 3566.95
$ws.Range("I126").Value = 1419.2
$ws.Range("K126").Value = 4257.6
$ws.Range("M126").Value = -1787.6

$ws.Range("H132").Value = 7701.675
$ws.Range("I132").Value = 8063.593
$ws.Range("J132").Value = 6950
$ws.Range("K132").Value = 24190.779
$ws.Range("L132").Value = 20850
$ws.Range("M132").Value = -21660.779
$ws.Range("N132").Value = -25910

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H43").Value = 26323
$ws.Range("I43").Value = 9000
$ws.Range("K43").Value = 9000
$ws.Range("M43").Value = -8851

$ws.Range("H49").Value = 20000
$ws.Range("I49").Value = 20000
$ws.Range("K49").Value = 20000
$ws.Range("M49").Value = -19770

$ws.Range("H113").Value = 1266.2222
$ws.Range("I113").Value = 863
$ws.Range("K113").Value = 2589
$ws.Range("M113").Value = -419

$ws.Range("H126").Value = 3846.077
$ws.Range("I126").Value = 3636.2727
$ws.Range("K126").Value = 10908.8181
$ws.Range("M126").Value = -8438.8181

$ws.Range("H131").Value = 80000
$ws.Range("J131").Value = 80000
$ws.Range("L131").Value = 80000
$ws.Range("N131").Value = -90080

Write-Host "Applied all Hyperion_Profits updates"
